$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# New data added to columns K/L/M for rows 34-38 (harmonogram related files)
# Use the same serial date number already used by sibling date cells (e.g. N34)
$newDate = 45887

$ws.Range("K34").Value = $newDate
$ws.Range("L34").Value = "WizytaControllerTests.cs"
$ws.Range("M34").Value = 10

$ws.Range("K35").Value = $newDate
$ws.Range("L35").Value = "HarmonogramControllerTest.cs"
$ws.Range("M35").Value = 48

$ws.Range("K36").Value = $newDate
$ws.Range("L36").Value = "HarmonogramService.cs"
$ws.Range("M36").Value = 1

$ws.Range("K37").Value = $newDate
$ws.Range("L37").Value = "HarmonogramController.cs"
$ws.Range("M37").Value = 7

$ws.Range("K38").Value = $newDate
$ws.Range("L38").Value = "MockWizytaRepository.cs"
$ws.Range("M38").Value = 4

# K39 just carries the date style but is left blank
# Copy formatting from an existing date cell (N34) so the same style index is reused
$ws.Range("N34").Copy()
$ws.Range("K34:K39").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the values that the formatting paste may have touched (safety)
$ws.Range("K34").Value = $newDate
$ws.Range("K35").Value = $newDate
$ws.Range("K36").Value = $newDate
$ws.Range("K37").Value = $newDate
$ws.Range("K38").Value = $newDate
$ws.Range("K39").ClearContents()

# Update selection to match the new active cell
$ws.Range("K39").Select()

$wb.Save()
